$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 16 (shifts old rows 16-18 down to 17-19,
# carrying their values/formatting with them).
$ws.Rows.Item(16).Insert()

# Refresh the "取得日時" (fetched-at) timestamp for every data row (2-19),
# matching the new run time recorded in the commit message.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = "2026-02-17 02:30:45"
}

# Populate the freshly inserted row 16 with the new listing.
$ws.Cells.Item(16, 2).Value = "【急募】よもぎ蒸しサロンのWebサイトエラー解決依頼"
$ws.Cells.Item(16, 3).Value = "システム開発"
$ws.Cells.Item(16, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(16, 5).Value = "期限情報なし"
$ws.Cells.Item(16, 6).Value = "https://www.lancers.jp/work/detail/5493140"
$ws.Cells.Item(16, 7).Value = 33
$ws.Cells.Item(16, 8).Value = "◇サイト"

# Row-insert doesn't renumber the existing hyperlink relationships, so rebuild
# the whole F-column hyperlink set from scratch in the correct top-to-bottom
# order (this also fixes up the one that used to live at the now-shifted rows).
$ws.Range("A1").Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5460562",
    "https://www.lancers.jp/work/detail/5473648",
    "https://www.lancers.jp/work/detail/5492832",
    "https://www.lancers.jp/work/detail/5217096",
    "https://www.lancers.jp/work/detail/5460563",
    "https://www.lancers.jp/work/detail/5488168",
    "https://www.lancers.jp/work/detail/5492887",
    "https://www.lancers.jp/work/detail/5492576",
    "https://www.lancers.jp/work/detail/5493016",
    "https://www.lancers.jp/work/detail/5492959",
    "https://www.lancers.jp/work/detail/5492441",
    "https://www.lancers.jp/work/detail/5492383",
    "https://www.lancers.jp/work/detail/5468432",
    "https://www.lancers.jp/work/detail/5492631",
    "https://www.lancers.jp/work/detail/5493140",
    "https://www.lancers.jp/work/detail/5492891",
    "https://www.lancers.jp/work/detail/5492894",
    "https://www.lancers.jp/work/detail/5492925"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($cell, $urls[$i])
    $cell.Style = "Hyperlink"
}
